# "edit table proxy and add functional for parser"
#
# The sheet previously had leftover/stale parsed rows (4-8, columns A:G)
# from a prior parser run: 362312/Peugeot---Citroen, 00004254A2/..,
# 00006426YN/.., 00008120T7/.., 6270000290/ГАЗ. Those are test rows left
# in the "table proxy" that the parser feeds - clear them out so only the
# two genuine sample rows (2 and 3) remain above the now-empty block that
# the parser will repopulate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the leftover parsed rows (4-8) across all their populated columns.
[void]$ws.Range("A4:G8").ClearContents()

# Move the selection onto the now-empty block ready for the parser to refill.
[void]$ws.Range("A4:G5").Select()
